$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 daily values - automatic electricity price update
$ws.Range("A2").Value = 45894
$ws.Range("B2").Value = 111.92
$ws.Range("C2").Value = 107.5
$ws.Range("D2").Value = 104.99
$ws.Range("E2").Value = 101.12
$ws.Range("F2").Value = 98.34999999999999
$ws.Range("G2").Value = 104.99
$ws.Range("H2").Value = 108.73
$ws.Range("I2").Value = 114.32
$ws.Range("J2").Value = 114.32
$ws.Range("K2").Value = 108.32
$ws.Range("L2").Value = 89.51000000000001
$ws.Range("M2").Value = 65.01000000000001
$ws.Range("N2").Value = 55.2
$ws.Range("O2").Value = 35
$ws.Range("P2").Value = 26.17
$ws.Range("Q2").Value = 25.2
$ws.Range("R2").Value = 56.43
$ws.Range("S2").Value = 70.09999999999999
$ws.Range("T2").Value = 97.43000000000001
$ws.Range("U2").Value = 114.78
$ws.Range("V2").Value = 125.95
$ws.Range("W2").Value = 142
$ws.Range("X2").Value = 123.11
$ws.Range("Y2").Value = 114.68
$ws.Range("Z2").Value = 92.3
$ws.Range("AB2").Value = 126.44
$ws.Range("AD2").Value = 133.98
$ws.Range("AF2").Value = 118.9
$ws.Range("AG2").Value = "10h-17h"
